$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:G2").Copy()
$ws.Range("A5:G5").PasteSpecial()

$ws.Range("A5").Value = "IC403"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "TQFP-128_14x14mm_Pitch0.4mm_EP"
$ws.Range("D5").Value = "XE216-512-TQ128"
$ws.Range("E5").Value = "xCore"
$ws.Range("F5").Value = "XMOS"
$ws.Range("G5").Value = "XE216-512-TQ128-C20"

$ws.Range("D20").Select()
